# Removed user from investor kyc
# The sheet originally had 11 columns:
#   A Investor | B Full Name | C First Name | D Last Name | E Email | F Phone |
#   G PAN | H Address | I Bank Account | J IFSC Code | K Send Confirmation Email
# The edit removes the per-user identity columns (First Name, Last Name,
# Email, Phone) and renames "Send Confirmation Email" to "Verified",
# normalizing both rows' value to "Yes".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the hyperlinks that lived on the (soon to be deleted) Email column
# before we shuffle columns around.
$ws.Hyperlinks.Delete()

# Remove columns C:F (First Name, Last Name, Email, Phone). Columns to the
# right shift left to fill the gap.
$ws.Range("C1:F1").EntireColumn.Delete()

# After the deletion the old "Send Confirmation Email" column (K) is now G.
# Rename the header and normalize both data rows to "Yes".
$ws.Range("G1").Value = "Verified"
$ws.Range("G2").Value = "Yes"
$ws.Range("G3").Value = "Yes"

# The hyperlink formatting is no longer used anywhere in the sheet; drop the
# now-orphaned "Hyperlink" cell style so it doesn't linger in the workbook.
$wb.Styles("Hyperlink").Delete()

# Match the saved cursor position from the authored workbook.
$ws.Range("G4").Select()
